## close #53 fix 2 monster big boss
## Insert a new skill row (55900045 "灭绝") into the Skill sheet, right
## before the "55990001" block (i.e. as the new row 176), pushing every
## subsequent row down by one and growing the "表3_25" table accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skill")

# --- 1. Insert a fresh row at 176, shifting 176:191 down to 177:192 ----
$ws.Rows.Item(176).Insert()

# Copy the formatting of the row above (175) onto the new row so the new
# row renders with the same style set used by its neighbours.
$ws.Range("A175:Z175").Copy()
$ws.Range("A176:Z176").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the taller row height used for the new skill's description.
$ws.Rows.Item(176).RowHeight = 84

# --- 2. Populate the new row's cells -----------------------------------
$ws.Range("A176").Value = 55900045
$ws.Range("B176").Value = "灭绝"
$ws.Range("C176").Value = "特殊"
$ws.Range("D176").Value = "NEA"
$ws.Range("E176").Value = 100
$ws.Range("F176").Value = "'true"
$ws.Range("G176").Value = ""
$ws.Range("H176").Value = "foreach(IMonster mon in s.Map.GetRangeMonster(s.IsLeft,sp.Target,sp.Shape,sp.Range,s.Position).FilterId(s.Id).SortDistance(true).Top(1)) mon.SuddenDeath();"
$ws.Range("Q176").Value = "Active"
$ws.Range("R176").Value = "'true"
$ws.Range("S176").Value = "召唤时杀死最近的敌人"
$ws.Range("X176").Value = 25
$ws.Range("Y176").Value = "miejue"

# --- 3. Grow the table / autofilter range to include the new row -------
$t = $ws.ListObjects.Item(1)
$t.Resize($ws.Range("A3:Z192"))

# --- 4. Leave the selection on the cell the author ended up on ---------
$ws.Range("R176").Select()
